# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed values, keeping the cells as plain text so
# figures like "46.546.96" or "0.999" are not reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text formatting before writing so Excel
# does not auto-coerce numeric-looking strings (e.g. "0.999", "9.79")
# into real numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "46.546.96"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "2.462.62"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "322.94"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "105.05"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "36.07"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "18.28"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "7.08"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "2.849.11"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "2.460.64"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "46.428.43"
$ws.Range("E18").Value = "  +4.25%  "
$ws.Range("D19").Value = "12.64"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "0.0₃0935"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "70.41"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").Value = "248.69"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("E24").Value = "  +4.26%  "
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").Value = "26.09"
$ws.Range("E26").Value = "  +3.35%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").Value = "9.79"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").Value = "35.37"
$ws.Range("E30").Value = "  +5.20%  "
$ws.Range("D31").Value = "49.46"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "19.57"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "5.33"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "0.0766"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("D38").Value = "1.90"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "122.87"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "20.65"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("D45").Value = "1.982.53"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("E48").Value = "  +3.73%  "
$ws.Range("D49").Value = "9.02"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").Value = "5.29"
$ws.Range("E50").Value = "  +15.04%  "
$ws.Range("D51").Value = "79.09"
$ws.Range("E51").Value = "  +4.92%  "

# Drop the temporary Text number format again so the cells end up
# back on the workbook default style (matching the original file).
$dataRange.ClearFormats()

